$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting the existing header row (A1:F1) one column to the right (-> B1:G1)
$ws.Columns.Item(1).Insert(-4161)
